$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ScreenRecStarted" was renamed to "0_unstated" throughout the transition-
# matrix labels (header + row labels that reference it).
$ws.Range("G1").Value = "0_unstated"
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Move the active selection to where the author last clicked.
$ws.Range("E14").Select()
